$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns for the rows being updated so that
# numeric-looking strings (e.g. "4.27") are stored as literal text, matching the
# original inline-string cell contents instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "58.745.07"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.493.61"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "535.67"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").Value = "136.62"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.514.73"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "0.346"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").Value = "2.936.30"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").Value = "22.89"
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "58.647.88"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "2.505.33"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "11.07"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("D20").Value = "4.27"
$ws.Range("D21").Value = "322.83"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D23").Value = "5.90"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "65.30"
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("D29").Value = "6.71"
$ws.Range("E29").Value = "  -3.41%  "
$ws.Range("D30").Value = "0.0₃0766"
$ws.Range("E30").Value = "  -1.92%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "167.06"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "18.40"
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("D37").Value = "4.08"
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("E38").Value = "  -3.54%  "
$ws.Range("D39").Value = "36.63"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("D40").Value = "0.813"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "3.59"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "285.43"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "131.20"
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("D46").Value = "0.603"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "0.0923"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("D50").Value = "0.0220"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").Value = "17.20"
$ws.Range("E51").Value = "  -3.79%  "
